$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1684397163120567
$ws.Range("C2").Value = 0.6223404255319149
$ws.Range("J2").Value = 0.01063829787234043
$ws.Range("P2").Value = 0.125886524822695
$ws.Range("S2").Value = 0.0726950354609929
$ws.Range("B3").Value = 0.008356545961002786
$ws.Range("C3").Value = 0.01671309192200557
$ws.Range("J3").Value = 0.02228412256267409
$ws.Range("P3").Value = 0.7103064066852368
$ws.Range("S3").Value = 0.2423398328690808
$ws.Range("J4").Value = 0.03846153846153846
$ws.Range("P4").Value = 0.7115384615384616
$ws.Range("S4").Value = 0.25
$ws.Range("B6").Value = 0.07755102040816327
$ws.Range("D6").Value = 0.01836734693877551
$ws.Range("F6").Value = 0.06734693877551021
$ws.Range("J6").Value = 0.2530612244897959
$ws.Range("O6").Value = 0.01224489795918367
$ws.Range("Q6").Value = 0.163265306122449
$ws.Range("R6").Value = 0.05714285714285714
$ws.Range("S6").Value = 0.3510204081632653
$ws.Range("B7").Value = 0.1157894736842105
$ws.Range("D7").Value = 0.02368421052631579
$ws.Range("E7").Value = 0.005263157894736842
$ws.Range("F7").Value = 0.06315789473684211
$ws.Range("J7").Value = 0.1210526315789474
$ws.Range("O7").Value = 0.01578947368421053
$ws.Range("Q7").Value = 0.1578947368421053
$ws.Range("R7").Value = 0.08157894736842106
$ws.Range("S7").Value = 0.4157894736842105
$ws.Range("B8").Value = 0.08695652173913043
$ws.Range("D8").Value = 0.01630434782608696
$ws.Range("E8").Value = 0.00108695652173913
$ws.Range("F8").Value = 0.07608695652173914
$ws.Range("J8").Value = 0.09456521739130434
$ws.Range("O8").Value = 0.0108695652173913
$ws.Range("Q8").Value = 0.1989130434782609
$ws.Range("R8").Value = 0.09130434782608696
$ws.Range("S8").Value = 0.4239130434782609
$ws.Range("B9").Value = 0.08514851485148515
$ws.Range("D9").Value = 0.01584158415841584
$ws.Range("E9").Value = 0.00198019801980198
$ws.Range("F9").Value = 0.05544554455445545
$ws.Range("J9").Value = 0.09504950495049505
$ws.Range("O9").Value = 0.01386138613861386
$ws.Range("Q9").Value = 0.2257425742574257
$ws.Range("R9").Value = 0.09108910891089109
$ws.Range("S9").Value = 0.4158415841584158
$ws.Range("B10").Value = 0.09773859716366425
$ws.Range("D10").Value = 0.02491376006132618
$ws.Range("F10").Value = 0.06822537370640092
$ws.Range("J10").Value = 0.121119202759678
$ws.Range("O10").Value = 0.01839785358374856
$ws.Range("Q10").Value = 0.2230739747029513
$ws.Range("R10").Value = 0.09658873131467996
$ws.Range("S10").Value = 0.3499425067075508
$ws.Range("G11").Value = 0.1563636363636364
$ws.Range("J11").Value = 0.08
$ws.Range("K11").Value = 0.1963636363636364
$ws.Range("L11").Value = 0.5600000000000001
$ws.Range("S11").Value = 0.007272727272727273
$ws.Range("G12").Value = 0.7563291139240507
$ws.Range("J12").Value = 0.180379746835443
$ws.Range("K12").Value = 0.006329113924050633
$ws.Range("L12").Value = 0.0189873417721519
$ws.Range("S12").Value = 0.0379746835443038
$ws.Range("G13").Value = 0.7045454545454546
$ws.Range("J13").Value = 0.1931818181818182
$ws.Range("S13").Value = 0.1022727272727273
$ws.Range("F15").Value = 0.02528735632183908
$ws.Range("H15").Value = 0.135632183908046
$ws.Range("I15").Value = 0.1057471264367816
$ws.Range("J15").Value = 0.367816091954023
$ws.Range("K15").Value = 0.05517241379310345
$ws.Range("M15").Value = 0.004597701149425287
$ws.Range("O15").Value = 0.0735632183908046
$ws.Range("S15").Value = 0.232183908045977
$ws.Range("F16").Value = 0.02077922077922078
$ws.Range("H16").Value = 0.1818181818181818
$ws.Range("I16").Value = 0.08311688311688312
$ws.Range("J16").Value = 0.3974025974025974
$ws.Range("K16").Value = 0.1402597402597403
$ws.Range("M16").Value = 0.007792207792207792
$ws.Range("N16").Value = 0.005194805194805195
$ws.Range("O16").Value = 0.06753246753246753
$ws.Range("S16").Value = 0.09610389610389611
$ws.Range("F17").Value = 0.01548886737657309
$ws.Range("H17").Value = 0.1771539206195547
$ws.Range("I17").Value = 0.1113262342691191
$ws.Range("J17").Value = 0.4181994191674734
$ws.Range("K17").Value = 0.06389157792836399
$ws.Range("M17").Value = 0.01452081316553727
$ws.Range("N17").Value = 0.000968054211035818
$ws.Range("O17").Value = 0.07163601161665054
$ws.Range("S17").Value = 0.1268151016456922
$ws.Range("F18").Value = 0.01769911504424779
$ws.Range("H18").Value = 0.1769911504424779
$ws.Range("I18").Value = 0.1349557522123894
$ws.Range("J18").Value = 0.3938053097345133
$ws.Range("K18").Value = 0.05973451327433629
$ws.Range("M18").Value = 0.01991150442477876
$ws.Range("O18").Value = 0.05752212389380531
$ws.Range("S18").Value = 0.1393805309734513
$ws.Range("F19").Value = 0.01788432267884323
$ws.Range("H19").Value = 0.2050989345509893
$ws.Range("I19").Value = 0.09817351598173515
$ws.Range("J19").Value = 0.3611111111111111
$ws.Range("K19").Value = 0.1008371385083714
$ws.Range("M19").Value = 0.02245053272450533
$ws.Range("N19").Value = 0.001902587519025875
$ws.Range("O19").Value = 0.05974124809741248
$ws.Range("S19").Value = 0.1328006088280061
